{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the finished requirement \"R12: Report the enemy that awards the\n// highest score and the level where it is located.\" and remove the whole\n// paragraph (it now lives implemented elsewhere, so the placeholder list\n// item goes away), leaving R7 and R13 as neighbours.\nconst target = paragraphs.items.find((p) => p.text.trim().indexOf(\"R12:\") === 0);\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the requirement paragraph \"R12: Report the enemy that awards the\n# highest score and the level where it is located.\" and remove it entirely\n# (including its paragraph mark), leaving R7 and R13 adjacent.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"R12:*highest score*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
